# Update countries & provincias Spain
#
# This script applies the data refresh captured in the commit:
#  - Two country rows swap places in the shared-string ordering; since the
#    row/column layout of the sheet stays fixed, that swap shows up here as
#    the "Camerun" row and "Estado de Palestina" row exchanging their
#    country labels (each row keeps its position, but now reports the other
#    country's figures, with Palestina's figures refreshed to the latest
#    count). Likewise for "Eslovenia" / "Lituania".
#  - A handful of country rows get refreshed case/recovery/death counts.
#  - The "last updated" timestamp banner in A1 is bumped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $country, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 1).Value = $country
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# Estados Unidos
Set-Row 4 "Estados Unidos" 5841689 261 3148080 2513427 0 8 180182

# Emiratos Arabes Unidos
Set-Row 45 "Emiratos Arabes Unidos" 67007 390 58488 8144 0 3 375

# Row 75 / 76 swap: Estado de Palestina <-> Camerun
Set-Row 75 "Estado de Palestina" 18802 326 11103 7572 0 2 127
Set-Row 76 "Camerun" 18762 0 16540 1814 0 0 408

# Row 129 / 130 swap: Lituania <-> Eslovenia
Set-Row 129 "Lituania" 2635 41 1766 785 0 0 84
Set-Row 130 "Eslovenia" 2617 0 2079 407 0 0 131

# Uganda
Set-Row 134 "Uganda" 2263 97 1199 1044 0 0 20

# Timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 23 de Agosto de 2020 a las 13:05"
